$d = $word.ActiveDocument

function Replace-Unique($findText, $replaceText) {
  $r = $d.Content
  $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
  if (-not $ok) { throw "Find/Replace failed for: $findText" }
}

# --- Title ---
Replace-Unique 'Unraveling the Enigmatic Universe' 'History: Unveiling the Tapestry of Human Endeavors'

# --- Author name: 'Sophia Robinson' -> 'Ms' + '.' + ' Olivia Alexander' (3 runs) ---
$p2 = $d.Paragraphs.Item(2).Range
$ok = $p2.Find.Execute('Sophia Robinson', $false, $false, $false, $false, $false, $true, 1, $false, 'Ms', 2)
if (-not $ok) { throw "Find/Replace failed for author name" }
$p2b = $d.Paragraphs.Item(2).Range
$okb = $p2b.Find.Execute('Ms')
if (-not $okb) { throw "Could not locate inserted Ms token" }
$p2b.Collapse(0)
$p2b.InsertAfter('.')
$p2c = $d.Paragraphs.Item(2).Range
$okc = $p2c.Find.Execute('Ms.')
if (-not $okc) { throw "Could not locate inserted Ms. token" }
$p2c.Collapse(0)
$p2c.InsertAfter(' Olivia Alexander')

# --- Email line: scoped to paragraph 3 to avoid ambiguous matches ---
$p3 = $d.Paragraphs.Item(3).Range
$ok = $p3.Find.Execute('sophia', $false, $false, $false, $false, $false, $true, 1, $false, 'olivia', 2)
if (-not $ok) { throw "Find/Replace failed for: sophia" }
$p3 = $d.Paragraphs.Item(3).Range
$ok = $p3.Find.Execute('robinson@xyz', $false, $false, $false, $false, $false, $true, 1, $false, 'alexander@schooledu', 2)
if (-not $ok) { throw "Find/Replace failed for: robinson@xyz" }
$p3 = $d.Paragraphs.Item(3).Range
$ok = $p3.Find.Execute('com', $false, $false, $false, $false, $false, $true, 1, $false, 'org', 2)
if (-not $ok) { throw "Find/Replace failed for: com -> org" }

# --- Body paragraph (paragraph 5): sentence-by-sentence replacements ---
Replace-Unique 'In the vast expanse of existence, humanity''s quest to comprehend the enigmas of the universe has been a relentless pursuit' 'History, like an intricate tapestry woven with threads of time, stands as a testament to the relentless march of human civilizations'
Replace-Unique ' From the celestial mechanics that orchestrate the cosmic ballet to the fundamental particles that underpin reality, our understanding of the universe has undergone a remarkable evolution' ' As we delve into the annals of the past, we embark on a journey through the triumphs and tribulations of our ancestors, gaining insights into the roots of our present and the seeds of our future'
Replace-Unique ' This journey of exploration has unveiled profound truths and elucidated perplexing mysteries, leaving us in awe of the intricate tapestry of cosmic existence' ' History offers a kaleidoscope of human experiences, from the grandeur of ancient empires to the complexities of modern societies, inviting us to ponder the actions and choices that have shaped our world'
Replace-Unique 'The symphony of celestial bodies, guided by the gravitational maestro, reveals intricate patterns and dynamic interactions' 'In the vast expanse of history, we encounter pivotal moments that have reshaped the course of human events'
Replace-Unique ' The interplay of stars, galaxies, and clusters, separated by unfathomable distances, captivates our imagination and invites us to unravel the mysteries of their formation and evolution' ' From the rise and fall of great civilizations to the transformative power of scientific discoveries, each era holds lessons that resonate with us today'
Replace-Unique ' As we delve deeper into the cosmos, we encounter cosmic phenomena that defy conventional understanding, such as black holes, wormholes, and dark matter, challenging our current scientific paradigms' ' History serves as a mirror, reflecting the strengths and weaknesses of humanity, offering poignant reminders of the consequences of our actions and the potential for progress when we embrace unity and understanding'
Replace-Unique 'The exploration of the subatomic realm has led to awe-inspiring discoveries, revealing a hidden world of particles and forces that govern the fundamental fabric of matter' 'History is not merely a chronicle of events; it is an exploration of the human condition, an ongoing dialogue between the past and the present'
Replace-Unique ' The Standard Model of Physics has illuminated the intricate dance of fundamental particles, providing a framework for understanding the forces that shape our universe' ' By studying history, we cultivate an appreciation for the diversity of human cultures, the resilience of the human spirit, and the interconnectedness of all things'
Replace-Unique ' Yet, tantalizing hints of undiscovered particles and forces continue to beckon us, inviting us to probe the deepest mysteries of the universe''s inner workings' ' History invites us to question our assumptions, to challenge conventional wisdom, and to seek out new perspectives, empowering us to become informed and engaged citizens in a rapidly changing world'

# --- Body paragraph (paragraph 5): append new content after the last replaced sentence ---
$p5 = $d.Paragraphs.Item(5).Range
$okp5 = $p5.Find.Execute('History invites us to question our assumptions, to challenge conventional wisdom, and to seek out new perspectives, empowering us to become informed and engaged citizens in a rapidly changing world')
if (-not $okp5) { throw "Could not locate body anchor sentence" }
$p5.Collapse(0)
$p5.InsertAfter('.' + [char]11 + [char]11 + 'Body:' + [char]11 + [char]11 + 'History encompasses a vast array of themes and subfields, each offering unique insights into the human experience. Political history chronicles the rise and fall of governments, the power struggles of leaders, and the impact of political decisions on the lives of ordinary people. Economic history examines the evolution of economic systems, trade routes, and the distribution of wealth, shedding light on the forces that drive economic growth and inequality. Social history explores the lives of ordinary people, their customs, beliefs, and everyday struggles, providing a glimpse into the fabric of societies past and present.' + [char]11 + [char]11 + 'Cultural history delves into the arts, literature, music, and traditions that define a people''s identity and heritage. Intellectual history traces the development of ideas, philosophies, and scientific advancements that have transformed our understanding of the world. Environmental history investigates the relationship between humans and their environment, highlighting the profound impact of human activities on the planet. By studying these diverse aspects of history, we gain a comprehensive understanding of the complexities of human societies and the factors that have shaped our world.' + [char]11 + [char]11 + 'History is not only about memorizing dates and facts; it is about developing critical thinking skills, analyzing evidence, and constructing informed arguments. By engaging in historical inquiry, students learn to evaluate sources, identify bias, and weigh competing interpretations. They develop the ability to think independently, to see connections between seemingly disparate events, and to make informed judgments about the past. These skills are essential for success in a wide range of fields and for navigating the challenges of the modern world')

# --- Summary heading is unchanged ---

# --- Summary paragraph (paragraph 7): sentence-by-sentence replacements ---
Replace-Unique 'Our exploration of the universe has been a symphony of wonder, revealing both profound truths and perplexing enigmas' 'History is an exploration of the human condition, a tapestry woven with the threads of time'
Replace-Unique ' From the cosmic ballet of celestial bodies to the subatomic world of fundamental particles, the quest for understanding has unlocked secrets and unveiled mysteries' ' Through the study of history, we gain insights into the triumphs and tribulations of our ancestors, the roots of our present, and the seeds of our future'
Replace-Unique ' The journey continues, and we stand at the threshold of new discoveries, poised to unravel the enigmas that still shroud the vast expanse of the universe' ' History encompasses a vast array of themes and subfields, from political and economic history to social, cultural, and intellectual history'

# --- Summary paragraph (paragraph 7): append new content ---
$p7 = $d.Paragraphs.Item(7).Range
$okp7 = $p7.Find.Execute('History encompasses a vast array of themes and subfields, from political and economic history to social, cultural, and intellectual history')
if (-not $okp7) { throw "Could not locate summary anchor sentence" }
$p7.Collapse(0)
$p7.InsertAfter('. By studying history, we cultivate an appreciation for the diversity of human cultures, the resilience of the human spirit, and the interconnectedness of all things. We develop critical thinking skills, learn to analyze evidence, and construct informed arguments, empowering us to become informed and engaged citizens in a rapidly changing world')

# --- Add a new empty paragraph at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Output $d.Content.Text
